$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the generic "Aplicação<US???>" placeholder text in column B (rows 14-22)
# with the specific user-story tags that were finally decided on.
$ws.Range("B14").Value = "Aplicação<US2.4>"
$ws.Range("B15").Value = "Aplicação<US2.4>"
$ws.Range("B16").Value = "Aplicação<US1.2>"
$ws.Range("B17").Value = "Aplicação<US2.3>"
$ws.Range("B18").Value = "Aplicação<US1.2>"
$ws.Range("B19").Value = "Aplicação<US2.3>"
$ws.Range("B20").Value = "Aplicação<US2.3>"
$ws.Range("B21").Value = "Aplicação<US1.2>"
$ws.Range("B22").Value = "Aplicação<US1.2>"

# Update the sheet selection/view to rest on the last edited cell.
$ws.Activate()
$ws.Range("B22").Select() | Out-Null
